$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 162, shifting existing rows 162-167 down to 163-168
$ws.Rows.Item(162).Insert()

# Populate the new row 162 with the new record
$ws.Cells.Item(162, 1).Value = 6
$ws.Cells.Item(162, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(162, 3).Value = "Metropolitana"
$ws.Cells.Item(162, 4).Value = 44516
$ws.Cells.Item(162, 5).Value = 13
$ws.Cells.Item(162, 6).Value = 100112022
$ws.Cells.Item(162, 7).Value = "Arveja Verde"
$ws.Cells.Item(162, 8).Value = "Sin especificar"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 350
$ws.Cells.Item(162, 11).Value = 13000
$ws.Cells.Item(162, 12).Value = 15000
$ws.Cells.Item(162, 13).Value = 14143
$ws.Cells.Item(162, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(162, 15).Value = "Región del Maule"
$ws.Cells.Item(162, 16).Value = 566
$ws.Cells.Item(162, 17).Value = 25
$ws.Cells.Item(162, 18).Value = "Hortaliza"
